$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26 ---
$ws.Range("Q26").Value2 = 639180
$ws.Range("R26").Value2 = 6701165
$ws.Range("Z26").ClearContents()
$ws.Range("AB26").ClearContents()

# --- Row 27 (becomes what used to be row 28's record) ---
$ws.Range("A27").Value2 = 112093192
$ws.Range("B27").Value2 = 90687
$ws.Range("D27").Value2 = "LC"
$ws.Range("E27").Value2 = 5964
$ws.Range("F27").Value2 = "Fjällig taggsvamp s.str."
$ws.Range("G27").Value2 = "Sarcodon imbricatus s.str."
$ws.Range("H27").Value2 = "(L.:Fr.) P.Karst."
$ws.Range("P27").Value2 = "Fagerdal, Upl"
$ws.Range("Q27").Value2 = 639180
$ws.Range("R27").Value2 = 6701165
$ws.Range("Z27").ClearContents()
$ws.Range("AB27").ClearContents()
$ws.Range("AC27").ClearContents()
$ws.Range("AD27").Value2 = $false

# --- Row 28 (becomes what used to be row 27's record) ---
$ws.Range("A28").Value2 = 112093171
$ws.Range("B28").Value2 = 88909
$ws.Range("D28").Value2 = "VU"
$ws.Range("E28").Value2 = 720
$ws.Range("F28").Value2 = "Violgubbe"
$ws.Range("G28").Value2 = "Gomphus clavatus"
$ws.Range("H28").Value2 = "(Pers.) Gray"
$ws.Range("P28").Value2 = "Bladsätra, Upl"
$ws.Range("Q28").Value2 = 639205
$ws.Range("R28").Value2 = 6701016
$ws.Range("Z28").ClearContents()
$ws.Range("AB28").ClearContents()
$ws.Range("AC28").Value2 = "Barkborredödat bestånd som även kantar mot ett stort kalhygge"
$ws.Range("AD28").Value2 = $true

# --- Row 29 ---
$ws.Range("Q29").Value2 = 639180
$ws.Range("R29").Value2 = 6701165
$ws.Range("Z29").ClearContents()
$ws.Range("AB29").ClearContents()
